$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range (rows 2..359 hold data; row 1 is the header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 359 }

# Column C ("Förändrad") moves from 45192 to 45202 for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}

# Row 3 ("A 6469-2022") also gets new tallies and two extra species names.
$ws.Cells.Item(3, 9).Value = 7    # I3 Signalarter
$ws.Cells.Item(3, 10).Value = 6   # J3 NT
$ws.Cells.Item(3, 15).Value = 9   # O3 Rödlistade
$ws.Cells.Item(3, 17).Value = 16  # Q3 Alla arter

$r3 = $ws.Cells.Item(3, 18)
$text = $r3.Value2
$text = $text -replace "Knärot`r`nGarnlav", "Knärot`r`nDofttaggsvamp`r`nGarnlav"
$text = $text -replace "Kornknutmossa`r`nVedticka", "Kornknutmossa`r`nStor revmossa`r`nVedticka"
$r3.Value = $text
